# Fruta / hortaliza, semanal
# Insert a new data row at row 2 (pushing all existing data rows down by one),
# matching the weekly update pattern: the newest record is prepended and the
# oldest historical rows shift down, with a brand-new row appearing at the
# bottom (row 28) that used to not exist.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right above the current row 2 (first data row), shifting
# all data down by one row.
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above (the bold/bordered
# header row). Clear that so it matches the look of a normal data row.
$ws.Rows.Item(2).ClearFormats()

# Column D holds dates and needs the same custom date format used by the
# rest of the data rows.
$ws.Cells.Item(2,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new first data row with the latest weekly record.
$ws.Cells.Item(2,1).Value  = 1
$ws.Cells.Item(2,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(2,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(2,4).Value  = 44498
$ws.Cells.Item(2,5).Value  = 15
$ws.Cells.Item(2,6).Value  = "Fruta"
$ws.Cells.Item(2,7).Value  = 100103
$ws.Cells.Item(2,8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(2,9).Value  = 100103004
$ws.Cells.Item(2,10).Value = "Durazno"
$ws.Cells.Item(2,11).Value = "Florida King"
$ws.Cells.Item(2,12).Value = "Segunda"
$ws.Cells.Item(2,13).Value = 200
$ws.Cells.Item(2,14).Value = 24000
$ws.Cells.Item(2,15).Value = 25000
$ws.Cells.Item(2,16).Value = 24500
$ws.Cells.Item(2,17).Value = "`$/bandeja 10 kilos granel"
$ws.Cells.Item(2,18).Value = "Región de Coquimbo"
$ws.Cells.Item(2,19).Value = 2450
$ws.Cells.Item(2,20).Value = 10
